$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.564.27"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "1.851.38"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.67"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6310"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9993"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07496"
$ws.Range("E8").Value = "  -0.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2915"
$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.78"
$ws.Range("E10").Value = "  +0.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07744"
$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("D12").Value = "1.852.12"
$ws.Range("E12").Value = "  +0.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.032"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6827"
$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001041"
$ws.Range("E15").Value = "  +0.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.63"
$ws.Range("E16").Value = "  -0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.269"
$ws.Range("E17").Value = "  +2.73%  "

$ws.Range("D18").Value = "29.558.45"
$ws.Range("E18").Value = "  +0.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "230.11"
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.41"
$ws.Range("E20").Value = "  +0.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.588"
$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9995"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.35"
$ws.Range("E24").Value = "  +0.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.539"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1373"
$ws.Range("E26").Value = "  -1.04%  "

$ws.Range("E27").Value = "  -0.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06712"
$ws.Range("E28").Value = "  +18.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.443"
$ws.Range("E29").Value = "  +0.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.487"
$ws.Range("E30").Value = "  +1.14%  "

$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("E32").Value = "  +1.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.842"
$ws.Range("E33").Value = "  +1.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.149"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6998"
$ws.Range("E35").Value = "  +0.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.576"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01868"
$ws.Range("E37").Value = "  +2.23%  "

$ws.Range("D38").Value = "1.267.76"
$ws.Range("E38").Value = "  +2.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.845"
$ws.Range("E39").Value = "  +4.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.850"
$ws.Range("E40").Value = "  +7.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9366"
$ws.Range("E41").Value = "  +4.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("D43").Value = "2.017.30"
$ws.Range("E43").Value = "  +0.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.30"
$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.35"
$ws.Range("E45").Value = "  +1.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.736"
$ws.Range("E46").Value = "  +3.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.117"
$ws.Range("E47").Value = "  -0.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1169"
$ws.Range("E48").Value = "  +1.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.052"
$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3965"
$ws.Range("E50").Value = "  -0.76%  "

$ws.Range("E51").Value = "  -3.67%  "
